$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - numeric "views" (F column) updates only
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitF = @{
    2  = 837
    6  = 1159
    8  = 56
    9  = 126
    11 = 1218
    15 = 894
    19 = 660
    21 = 1750
    22 = 3126
    23 = 913
    25 = 2296
    27 = 8
    28 = 3151
    29 = 647
    30 = 745
    36 = 40
    37 = 104
    38 = 1121
    39 = 1812
    42 = 562
}
foreach ($row in $exhibitF.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitF[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - numeric "views" (F column) update
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(12, 6).Value = 93

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - numeric "views" (F column) updates
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allF = @{
    2  = 837
    5  = 1159
    7  = 126
    8  = 1218
    11 = 894
    18 = 1750
    19 = 3126
    20 = 913
    23 = 2296
    25 = 3151
    26 = 647
    27 = 745
    34 = 93
    41 = 1121
    42 = 1812
    45 = 562
}
foreach ($row in $allF.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allF[$row]
}

# Row 24 on "全部类型": old event replaced by a new one
$wsAll.Cells.Item(24, 3).Value = "杭州·天空漫境-第十二届（免费展）"
$wsAll.Cells.Item(24, 4).Value = "金桥北路990号 万达广场(杭州富阳店)"
$wsAll.Cells.Item(24, 5).Value = "2024.05.02 10:00-05.02 16:00"
$wsAll.Cells.Item(24, 6).Value = 8
$wsAll.Cells.Item(24, 7).Value = 25
$wsAll.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84057"
$wsAll.Cells.Item(24, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/g9uui75m1712574749479.jpeg"

# Row 38 on "全部类型": new event inserted in place (values overwritten)
$wsAll.Cells.Item(38, 2).NumberFormat = "@"
$wsAll.Cells.Item(38, 2).Value = "2024-05-26"
$wsAll.Cells.Item(38, 3).Value = "杭州·恋与深空×恋与制作人only"
$wsAll.Cells.Item(38, 4).Value = "望江东路333号 杭州瑞莱克斯大酒店"
$wsAll.Cells.Item(38, 5).Value = "2024.05.26 10:00-05.26 17:00"
$wsAll.Cells.Item(38, 6).Value = 40
$wsAll.Cells.Item(38, 7).Value = 60
$wsAll.Cells.Item(38, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84077"
$wsAll.Cells.Item(38, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/V6V4Pppv1712736555042.jpeg"

# Row 39 on "全部类型": now holds what used to be row 38's event (with updated view count)
$wsAll.Cells.Item(39, 2).NumberFormat = "@"
$wsAll.Cells.Item(39, 2).Value = "2024-06-01"
$wsAll.Cells.Item(39, 3).Value = "杭州·造梦探险家——二次元同好会"
$wsAll.Cells.Item(39, 4).Value = "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
$wsAll.Cells.Item(39, 5).Value = "2024.06.01 10:00-06.01 16:00"
$wsAll.Cells.Item(39, 6).Value = 104
$wsAll.Cells.Item(39, 7).Value = 28
$wsAll.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82736"
$wsAll.Cells.Item(39, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/lqXD63661711623533572.png"
